$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157; this shifts rows 157-255 down to 158-256.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new data record.
$ws.Cells.Item(157, 1).Value = 1
$ws.Cells.Item(157, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(157, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(157, 4).Value = 44726
$ws.Cells.Item(157, 5).Value = 15
$ws.Cells.Item(157, 6).Value = "Fruta"
$ws.Cells.Item(157, 7).Value = 100108
$ws.Cells.Item(157, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(157, 9).Value = 100108006
$ws.Cells.Item(157, 10).Value = "Plátano"
$ws.Cells.Item(157, 11).Value = "Barraganete"
$ws.Cells.Item(157, 12).Value = "Primera"
$ws.Cells.Item(157, 13).Value = 120
$ws.Cells.Item(157, 14).Value = 21000
$ws.Cells.Item(157, 15).Value = 22000
$ws.Cells.Item(157, 16).Value = 21500
$ws.Cells.Item(157, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(157, 18).Value = "Ecuador"
$ws.Cells.Item(157, 19).Value = 1075
$ws.Cells.Item(157, 20).Value = 20
